$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("Summary")
$wsCards = $wb.Worksheets.Item("Cards_telegram")

# Summary row 2
$wsSummary.Range("A2").Value = 1369
$wsSummary.Range("B2").Value = 45987.52083333334
$wsSummary.Range("C2").Value = "Адмирал"
$wsSummary.Range("D2").Value = "Амур"
$wsSummary.Range("E2").Value = "Адмирал – Амур"
$wsSummary.Range("F2").Value = 897818
$wsSummary.Range("G2").Value = "https://text.khl.ru/text/897818.html"
$wsSummary.Range("H2").Value = 1.91778
$wsSummary.Range("I2").Value = 2.626433
$wsSummary.Range("J2").Value = 4.544213
$wsSummary.Range("K2").Value = 29.99177
$wsSummary.Range("L2").Value = 28.61054
$wsSummary.Range("M2").Value = 58.60231
$wsSummary.Range("N2").Value = 0.625932
$wsSummary.Range("O2").Value = 0.153208
$wsSummary.Range("P2").Value = 0.220113
$wsSummary.Range("Q2").Value = 1.597617632586287
$wsSummary.Range("R2").Value = 6.527074304213879
$wsSummary.Range("S2").Value = 4.543121033287448
$wsSummary.Range("T2").Value = 62.5932
$wsSummary.Range("U2").Value = 15.3208
$wsSummary.Range("V2").Value = 22.0113
$wsSummary.Range("W2").Value = 0.365474
$wsSummary.Range("X2").Value = 0.633778
$wsSummary.Range("Y2").Value = 1.577839558962287
$wsSummary.Range("Z2").Value = 0.5376300000000001
$wsSummary.Range("AA2").Value = 0.461622
$wsSummary.Range("AB2").Value = 2.166274570969321
$wsSummary.Range("AC2").Value = 0.693985
$wsSummary.Range("AD2").Value = 0.305268
$wsSummary.Range("AE2").Value = 3.275810107839669
$wsSummary.Range("AF2").Value = 0.848425
$wsSummary.Range("AG2").Value = 0.151575
$wsSummary.Range("AH2").Value = 1.178654565813124
$wsSummary.Range("AI2").Value = 0.652263
$wsSummary.Range("AJ2").Value = 0.347737
$wsSummary.Range("AK2").Value = 1.533123908607418
$wsSummary.Range("AL2").Value = 0.61791
$wsSummary.Range("AM2").Value = 0.38209
$wsSummary.Range("AN2").Value = 1.618358660646372
$wsSummary.Range("AO2").Value = 0.347775
$wsSummary.Range("AP2").Value = 0.6522250000000001
$wsSummary.Range("AQ2").Value = 2.875422327654374
$wsSummary.Range("AR2").Value = 0.8879050000000001
$wsSummary.Range("AS2").Value = 1.126246614221116
$wsSummary.Range("AT2").Value = 0.548087
$wsSummary.Range("AU2").Value = 1.824527857803597

# Summary row 3
$wsSummary.Range("A3").Value = 1369
$wsSummary.Range("B3").Value = 45987.70833333334
$wsSummary.Range("C3").Value = "Трактор"
$wsSummary.Range("D3").Value = "Драконы"
$wsSummary.Range("E3").Value = "Трактор – Драконы"
$wsSummary.Range("F3").Value = 897816
$wsSummary.Range("G3").Value = "https://text.khl.ru/text/897816.html"
$wsSummary.Range("H3").Value = 4.435928
$wsSummary.Range("I3").Value = 3.568607
$wsSummary.Range("J3").Value = 8.004535000000001
$wsSummary.Range("K3").Value = 39.059641
$wsSummary.Range("L3").Value = 31.713169
$wsSummary.Range("M3").Value = 70.77281000000001
$wsSummary.Range("N3").Value = 0.300991
$wsSummary.Range("O3").Value = 0.139072
$wsSummary.Range("P3").Value = 0.554329
$wsSummary.Range("Q3").Value = 3.322358475834826
$wsSummary.Range("R3").Value = 7.190520018407732
$wsSummary.Range("S3").Value = 1.803982833299358
$wsSummary.Range("T3").Value = 30.0991
$wsSummary.Range("U3").Value = 13.9072
$wsSummary.Range("V3").Value = 55.4329
$wsSummary.Range("W3").Value = 0.11946
$wsSummary.Range("X3").Value = 0.874932
$wsSummary.Range("Y3").Value = 1.142945966086507
$wsSummary.Range("Z3").Value = 0.22231
$wsSummary.Range("AA3").Value = 0.772082
$wsSummary.Range("AB3").Value = 1.295199214591196
$wsSummary.Range("AC3").Value = 0.353962
$wsSummary.Range("AD3").Value = 0.6404300000000001
$wsSummary.Range("AE3").Value = 1.561450900176444
$wsSummary.Range("AF3").Value = 0.850654
$wsSummary.Range("AG3").Value = 0.149346
$wsSummary.Range("AH3").Value = 1.175566093852495
$wsSummary.Range("AI3").Value = 0.656018
$wsSummary.Range("AJ3").Value = 0.343982
$wsSummary.Range("AK3").Value = 1.524348417269038
$wsSummary.Range("AL3").Value = 0.928213
$wsSummary.Range("AM3").Value = 0.071787
$wsSummary.Range("AN3").Value = 1.07733892974996
$wsSummary.Range("AO3").Value = 0.802916
$wsSummary.Range("AP3").Value = 0.197084
$wsSummary.Range("AQ3").Value = 1.245460297216645
$wsSummary.Range("AR3").Value = 0.5863159999999999
$wsSummary.Range("AS3").Value = 1.705564917211879
$wsSummary.Range("AT3").Value = 0.80824
$wsSummary.Range("AU3").Value = 1.237256260516678

# Summary row 4
$wsSummary.Range("A4").Value = 1369
$wsSummary.Range("B4").Value = 45987.79166666666
$wsSummary.Range("C4").Value = "Северсталь"
$wsSummary.Range("D4").Value = "СКА"
$wsSummary.Range("E4").Value = "Северсталь – СКА"
$wsSummary.Range("F4").Value = 897817
$wsSummary.Range("G4").Value = "https://text.khl.ru/text/897817.html"
$wsSummary.Range("H4").Value = 1.464286
$wsSummary.Range("I4").Value = 2.390395
$wsSummary.Range("J4").Value = 3.854681
$wsSummary.Range("K4").Value = 24.993415
$wsSummary.Range("L4").Value = 27.964901
$wsSummary.Range("M4").Value = 52.958316
$wsSummary.Range("N4").Value = 0.317655
$wsSummary.Range("O4").Value = 0.231914
$wsSummary.Range("P4").Value = 0.450429
$wsSummary.Range("Q4").Value = 3.148069446411988
$wsSummary.Range("R4").Value = 4.311943220331674
$wsSummary.Range("S4").Value = 2.220105721434455
$wsSummary.Range("T4").Value = 31.7655
$wsSummary.Range("U4").Value = 23.1914
$wsSummary.Range("V4").Value = 45.0429
$wsSummary.Range("W4").Value = 0.783058
$wsSummary.Range("X4").Value = 0.216939
$wsSummary.Range("Y4").Value = 4.609590714440465
$wsSummary.Range("Z4").Value = 0.896154
$wsSummary.Range("AA4").Value = 0.103844
$wsSummary.Range("AB4").Value = 9.62982935942375
$wsSummary.Range("AC4").Value = 0.956212
$wsSummary.Range("AD4").Value = 0.043785
$wsSummary.Range("AE4").Value = 22.83887175973507
$wsSummary.Range("AF4").Value = 0.421034
$wsSummary.Range("AG4").Value = 0.578966
$wsSummary.Range("AH4").Value = 2.375105098400604
$wsSummary.Range("AI4").Value = 0.17563
$wsSummary.Range("AJ4").Value = 0.82437
$wsSummary.Range("AK4").Value = 5.693788077207766
$wsSummary.Range("AL4").Value = 0.521747
$wsSummary.Range("AM4").Value = 0.478253
$wsSummary.Range("AN4").Value = 1.916637757380493
$wsSummary.Range("AO4").Value = 0.255703
$wsSummary.Range("AP4").Value = 0.744297
$wsSummary.Range("AQ4").Value = 3.910787124124472
$wsSummary.Range("AR4").Value = 0.7598
$wsSummary.Range("AS4").Value = 1.316135825217162
$wsSummary.Range("AT4").Value = 0.855151
$wsSummary.Range("AU4").Value = 1.169384120465275

# Summary row 5
$wsSummary.Range("A5").Value = 1369
$wsSummary.Range("B5").Value = 45987.8125
$wsSummary.Range("B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsSummary.Range("C5").Value = "Динамо М"
$wsSummary.Range("D5").Value = "Локомотив"
$wsSummary.Range("E5").Value = "Динамо М – Локомотив"
$wsSummary.Range("F5").Value = 897815
$wsSummary.Range("G5").Value = "https://text.khl.ru/text/897815.html"
$wsSummary.Range("H5").Value = 2.014963
$wsSummary.Range("I5").Value = 2.067992
$wsSummary.Range("J5").Value = 4.082955
$wsSummary.Range("K5").Value = 24.33595
$wsSummary.Range("L5").Value = 26.124622
$wsSummary.Range("M5").Value = 50.460572
$wsSummary.Range("N5").Value = 0.311642
$wsSummary.Range("O5").Value = 0.194412
$wsSummary.Range("P5").Value = 0.493899
$wsSummary.Range("Q5").Value = 3.208810109035368
$wsSummary.Range("R5").Value = 5.143715408513878
$wsSummary.Range("S5").Value = 2.024705455973792
$wsSummary.Range("T5").Value = 31.1642
$wsSummary.Range("U5").Value = 19.4412
$wsSummary.Range("V5").Value = 49.3899
$wsSummary.Range("W5").Value = 0.572371
$wsSummary.Range("X5").Value = 0.427582
$wsSummary.Range("Y5").Value = 2.338732687531281
$wsSummary.Range("Z5").Value = 0.738323
$wsSummary.Range("AA5").Value = 0.261631
$wsSummary.Range("AB5").Value = 3.822177035595935
$wsSummary.Range("AC5").Value = 0.856978
$wsSummary.Range("AD5").Value = 0.142975
$wsSummary.Range("AE5").Value = 6.994229760447631
$wsSummary.Range("AF5").Value = 0.565997
$wsSummary.Range("AG5").Value = 0.434003
$wsSummary.Range("AH5").Value = 1.766793816928358
$wsSummary.Range("AI5").Value = 0.296038
$wsSummary.Range("AJ5").Value = 0.703962
$wsSummary.Range("AK5").Value = 3.377944723312547
$wsSummary.Range("AL5").Value = 0.689572
$wsSummary.Range("AM5").Value = 0.310428
$wsSummary.Range("AN5").Value = 1.450174891091866
$wsSummary.Range("AO5").Value = 0.427911
$wsSummary.Range("AP5").Value = 0.572089
$wsSummary.Range("AQ5").Value = 2.336934549474073
$wsSummary.Range("AR5").Value = 0.6964630000000001
$wsSummary.Range("AS5").Value = 1.435826454528094
$wsSummary.Range("AT5").Value = 0.839553
$wsSummary.Range("AU5").Value = 1.19111003117135

# Cards_telegram row 2
$wsCards.Range("A2").Value = 45987.52083333334
$wsCards.Range("B2").Value = "Адмирал – Амур"
$text_2 = @"
КХЛ • Регулярный чемпионат • 26.11.2025
Адмирал – Амур
Ожидания модели (60’):
• Голы: λ_total ≈ 5.45 (3.36 : 2.09)
• Броски: SOG λ ≈ 59 (30 : 29)
Исход (60’), честные кф:
• П1: 62.6%  (Kмод 1.60)
• Х:  15.3%  (Kмод 6.53)
• П2: 22.0%  (Kмод 4.54)
Тоталы голов:
• ТМ 4.5: 36.5%  (Kмод 2.74)
• ТБ 4.5: 63.4%  (Kмод 1.58)
• ТМ 5.5: 53.8%  (Kмод 1.86)
• ТБ 5.5: 46.2%  (Kмод 2.17)
• ТМ 6.5: 69.4%  (Kмод 1.44)
• ТБ 6.5: 30.5%  (Kмод 3.28)
Индивидуальные тоталы:
• Адмирал ИТБ 1.5: 84.8% (Kмод 1.18)
• Адмирал ИТБ 2.5: 65.2% (Kмод 1.53)
• Амур ИТБ 1.5: 61.8% (Kмод 1.62)
• Амур ИТБ 2.5: 34.8% (Kмод 2.88)
Фора +1.5:
• Адмирал +1.5: 88.8% (Kмод 1.13)
• Амур +1.5: 54.8% (Kмод 1.82)
"@
$wsCards.Range("C2").Value = $text_2

# Cards_telegram row 3
$wsCards.Range("A3").Value = 45987.70833333334
$wsCards.Range("B3").Value = "Трактор – Драконы"
$text_3 = @"
КХЛ • Регулярный чемпионат • 26.11.2025
Трактор – Драконы
Ожидания модели (60’):
• Голы: λ_total ≈ 7.68 (3.38 : 4.30)
• Броски: SOG λ ≈ 71 (39 : 32)
Исход (60’), честные кф:
• П1: 30.1%  (Kмод 3.32)
• Х:  13.9%  (Kмод 7.19)
• П2: 55.4%  (Kмод 1.80)
Тоталы голов:
• ТМ 4.5: 11.9%  (Kмод 8.37)
• ТБ 4.5: 87.5%  (Kмод 1.14)
• ТМ 5.5: 22.2%  (Kмод 4.50)
• ТБ 5.5: 77.2%  (Kмод 1.30)
• ТМ 6.5: 35.4%  (Kмод 2.83)
• ТБ 6.5: 64.0%  (Kмод 1.56)
Индивидуальные тоталы:
• Трактор ИТБ 1.5: 85.1% (Kмод 1.18)
• Трактор ИТБ 2.5: 65.6% (Kмод 1.52)
• Драконы ИТБ 1.5: 92.8% (Kмод 1.08)
• Драконы ИТБ 2.5: 80.3% (Kмод 1.25)
Фора +1.5:
• Трактор +1.5: 58.6% (Kмод 1.71)
• Драконы +1.5: 80.8% (Kмод 1.24)
"@
$wsCards.Range("C3").Value = $text_3

# Cards_telegram row 4
$wsCards.Range("A4").Value = 45987.79166666666
$wsCards.Range("B4").Value = "Северсталь – СКА"
$text_4 = @"
КХЛ • Регулярный чемпионат • 26.11.2025
Северсталь – СКА
Ожидания модели (60’):
• Голы: λ_total ≈ 3.19 (1.44 : 1.75)
• Броски: SOG λ ≈ 53 (25 : 28)
Исход (60’), честные кф:
• П1: 31.8%  (Kмод 3.15)
• Х:  23.2%  (Kмод 4.31)
• П2: 45.0%  (Kмод 2.22)
Тоталы голов:
• ТМ 4.5: 78.3%  (Kмод 1.28)
• ТБ 4.5: 21.7%  (Kмод 4.61)
• ТМ 5.5: 89.6%  (Kмод 1.12)
• ТБ 5.5: 10.4%  (Kмод 9.63)
• ТМ 6.5: 95.6%  (Kмод 1.05)
• ТБ 6.5: 4.4%  (Kмод 22.84)
Индивидуальные тоталы:
• Северсталь ИТБ 1.5: 42.1% (Kмод 2.38)
• Северсталь ИТБ 2.5: 17.6% (Kмод 5.69)
• СКА ИТБ 1.5: 52.2% (Kмод 1.92)
• СКА ИТБ 2.5: 25.6% (Kмод 3.91)
Фора +1.5:
• Северсталь +1.5: 76.0% (Kмод 1.32)
• СКА +1.5: 85.5% (Kмод 1.17)
"@
$wsCards.Range("C4").Value = $text_4

# Cards_telegram row 5
$wsCards.Range("A5").Value = 45987.8125
$wsCards.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsCards.Range("B5").Value = "Динамо М – Локомотив"
$text_5 = @"
КХЛ • Регулярный чемпионат • 26.11.2025
Динамо М – Локомотив
Ожидания модели (60’):
• Голы: λ_total ≈ 4.29 (1.90 : 2.39)
• Броски: SOG λ ≈ 50 (24 : 26)
Исход (60’), честные кф:
• П1: 31.2%  (Kмод 3.21)
• Х:  19.4%  (Kмод 5.14)
• П2: 49.4%  (Kмод 2.02)
Тоталы голов:
• ТМ 4.5: 57.2%  (Kмод 1.75)
• ТБ 4.5: 42.8%  (Kмод 2.34)
• ТМ 5.5: 73.8%  (Kмод 1.35)
• ТБ 5.5: 26.2%  (Kмод 3.82)
• ТМ 6.5: 85.7%  (Kмод 1.17)
• ТБ 6.5: 14.3%  (Kмод 6.99)
Индивидуальные тоталы:
• Динамо М ИТБ 1.5: 56.6% (Kмод 1.77)
• Динамо М ИТБ 2.5: 29.6% (Kмод 3.38)
• Локомотив ИТБ 1.5: 69.0% (Kмод 1.45)
• Локомотив ИТБ 2.5: 42.8% (Kмод 2.34)
Фора +1.5:
• Динамо М +1.5: 69.6% (Kмод 1.44)
• Локомотив +1.5: 84.0% (Kмод 1.19)
"@
$wsCards.Range("C5").Value = $text_5
